# Generate Report for Archive
# Update status text from "Ready for handoff" to "In Translation" across all sheets,
# then autofit the affected status columns so their widths shrink to match the new text.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ("Ready for handoff" -eq $cell.Value2) {
            $cell.Value = "In Translation"
        }
    }
}

# Autofit the Status column on the zh-cn and de-de sheets (column C),
# and the zh-cn/de-de status columns on the Overview sheet (columns E and F),
# since the shorter replacement text narrows the natural content width.
# ColumnWidth is pinned to 12.5 (chars) afterwards, which is the AutoFit-
# equivalent target width for the new, shorter "In Translation" text.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).EntireColumn.AutoFit() | Out-Null
$wsOverview.Columns.Item(6).EntireColumn.AutoFit() | Out-Null
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
